$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 to F1, then set header text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in time_taken values for each data row
$timeTaken = @(
    "2021-10-05 13:39:31.326783",
    "2021-10-05 13:39:31.326793",
    "2021-10-05 13:39:31.326796",
    "2021-10-05 13:39:31.326799",
    "2021-10-05 13:39:31.326802",
    "2021-10-05 13:39:31.326804",
    "2021-10-05 13:39:31.326807",
    "2021-10-05 13:39:31.326809",
    "2021-10-05 13:39:31.326812",
    "2021-10-05 13:39:31.326815",
    "2021-10-05 13:39:31.326817",
    "2021-10-05 13:39:31.326820",
    "2021-10-05 13:39:31.326822",
    "2021-10-05 13:39:31.326825",
    "2021-10-05 13:39:31.326827",
    "2021-10-05 13:39:31.326830",
    "2021-10-05 13:39:31.326832",
    "2021-10-05 13:39:31.326835",
    "2021-10-05 13:39:31.326838",
    "2021-10-05 13:39:31.326840",
    "2021-10-05 13:39:31.326843",
    "2021-10-05 13:39:31.326845",
    "2021-10-05 13:39:31.326848",
    "2021-10-05 13:39:31.326850",
    "2021-10-05 13:39:31.326853",
    "2021-10-05 13:39:31.326855",
    "2021-10-05 13:39:31.326858",
    "2021-10-05 13:39:31.326860",
    "2021-10-05 13:39:31.326863",
    "2021-10-05 13:39:31.326865",
    "2021-10-05 13:39:31.326867",
    "2021-10-05 13:39:31.326870",
    "2021-10-05 13:39:31.326873",
    "2021-10-05 13:39:31.326875",
    "2021-10-05 13:39:31.326878",
    "2021-10-05 13:39:31.326880",
    "2021-10-05 13:39:31.326882",
    "2021-10-05 13:39:31.326885",
    "2021-10-05 13:39:31.326887",
    "2021-10-05 13:39:31.326890",
    "2021-10-05 13:39:31.326893",
    "2021-10-05 13:39:31.326895",
    "2021-10-05 13:39:31.326898",
    "2021-10-05 13:39:31.326900",
    "2021-10-05 13:39:31.326903",
    "2021-10-05 13:39:31.326905",
    "2021-10-05 13:39:31.326907",
    "2021-10-05 13:39:31.326910",
    "2021-10-05 13:39:31.326912",
    "2021-10-05 13:39:31.326915",
    "2021-10-05 13:39:31.326917",
    "2021-10-05 13:39:31.326919",
    "2021-10-05 13:39:31.326922",
    "2021-10-05 13:39:31.326925",
    "2021-10-05 13:39:31.326927",
    "2021-10-05 13:39:31.326930",
    "2021-10-05 13:39:31.326932",
    "2021-10-05 13:39:31.326935",
    "2021-10-05 13:39:31.326937",
    "2021-10-05 13:39:31.326940",
    "2021-10-05 13:39:31.326942",
    "2021-10-05 13:39:31.326944",
    "2021-10-05 13:39:31.326947",
    "2021-10-05 13:39:31.326949",
    "2021-10-05 13:39:31.326953",
    "2021-10-05 13:39:31.326956",
    "2021-10-05 13:39:31.326958",
    "2021-10-05 13:39:31.326960",
    "2021-10-05 13:39:31.326963",
    "2021-10-05 13:39:31.326965",
    "2021-10-05 13:39:31.326968",
    "2021-10-05 13:39:31.326970",
    "2021-10-05 13:39:31.326973",
    "2021-10-05 13:39:31.326975",
    "2021-10-05 13:39:31.326978",
    "2021-10-05 13:39:31.326980",
    "2021-10-05 13:39:31.326984",
    "2021-10-05 13:39:31.326987",
    "2021-10-05 13:39:31.326990",
    "2021-10-05 13:39:31.326992",
    "2021-10-05 13:39:31.326995",
    "2021-10-05 13:39:31.326997",
    "2021-10-05 13:39:31.327000",
    "2021-10-05 13:39:31.327002",
    "2021-10-05 13:39:31.327005",
    "2021-10-05 13:39:31.327007",
    "2021-10-05 13:39:31.327010",
    "2021-10-05 13:39:31.327012",
    "2021-10-05 13:39:31.327014",
    "2021-10-05 13:39:31.327017",
    "2021-10-05 13:39:31.327019",
    "2021-10-05 13:39:31.327022",
    "2021-10-05 13:39:31.327025",
    "2021-10-05 13:39:31.327028",
    "2021-10-05 13:39:31.327031",
    "2021-10-05 13:39:31.327033",
    "2021-10-05 13:39:31.327035",
    "2021-10-05 13:39:31.327038",
    "2021-10-05 13:39:31.327040",
    "2021-10-05 13:39:31.327043",
    "2021-10-05 13:39:31.327045",
    "2021-10-05 13:39:31.327047",
    "2021-10-05 13:39:31.327050",
    "2021-10-05 13:39:31.327052",
    "2021-10-05 13:39:31.327055",
    "2021-10-05 13:39:31.327057",
    "2021-10-05 13:39:31.327060",
    "2021-10-05 13:39:31.327062",
    "2021-10-05 13:39:31.327066",
    "2021-10-05 13:39:31.327069",
    "2021-10-05 13:39:31.327072",
    "2021-10-05 13:39:31.327074",
    "2021-10-05 13:39:31.327077",
    "2021-10-05 13:39:31.327079",
    "2021-10-05 13:39:31.327082",
    "2021-10-05 13:39:31.327084",
    "2021-10-05 13:39:31.327086",
    "2021-10-05 13:39:31.327089",
    "2021-10-05 13:39:31.327091",
    "2021-10-05 13:39:31.327094",
    "2021-10-05 13:39:31.327096",
    "2021-10-05 13:39:31.327099",
    "2021-10-05 13:39:31.327101",
    "2021-10-05 13:39:31.327103",
    "2021-10-05 13:39:31.327106",
    "2021-10-05 13:39:31.327108",
    "2021-10-05 13:39:31.327111",
    "2021-10-05 13:39:31.327113",
    "2021-10-05 13:39:31.327117",
    "2021-10-05 13:39:31.327120",
    "2021-10-05 13:39:31.327123",
    "2021-10-05 13:39:31.327126",
    "2021-10-05 13:39:31.327128",
    "2021-10-05 13:39:31.327131",
    "2021-10-05 13:39:31.327133",
    "2021-10-05 13:39:31.327136",
    "2021-10-05 13:39:31.327138",
    "2021-10-05 13:39:31.327141",
    "2021-10-05 13:39:31.327143",
    "2021-10-05 13:39:31.327146",
    "2021-10-05 13:39:31.327148",
    "2021-10-05 13:39:31.327151",
    "2021-10-05 13:39:31.327154",
    "2021-10-05 13:39:31.327156",
    "2021-10-05 13:39:31.327159",
    "2021-10-05 13:39:31.327161",
    "2021-10-05 13:39:31.327164",
    "2021-10-05 13:39:31.327166",
    "2021-10-05 13:39:31.327169",
    "2021-10-05 13:39:31.327172"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

